# Upload timesheet 26/6/18
# Fill in the timesheet row for Wednesday 26 June 2018 (row 5), which
# was previously blank.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date column (A) keeps its existing date style; just give it the date.
$ws.Range("A5").Value = 43277   # 26 June 2018

# Time columns (B:E) - match the h:mm time format already used by the
# rows above (this reuses the existing style rather than creating a new one).
$ws.Range("B5:E5").NumberFormat = "h:mm"
$ws.Range("B5").Value = 0.375                 # 9:00 - clock in (Nic)
$ws.Range("C5").Value = 0.1875                # 4:30 - clock out (Nic)
$ws.Range("D5").Value = 0.375                 # 9:00 - clock in (Arpit)
$ws.Range("E5").Value = 0.20833333333333334   # 5:00 - clock out (Arpit)

# Move the active selection to F11, matching where the editor left off.
[void]$ws.Range("F11").Select()
